$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.750.38'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '3.404.20'
$ws.Range('E3').Value = '  -1.69%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.92'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.12'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.89%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '3.405.75'
$ws.Range('E8').Value = '  -1.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.569'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -6.57%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.16'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.119'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -5.02%  '
$ws.Range('E12').Value = '  -5.59%  '
$ws.Range('D13').Value = '3.996.13'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('E14').Value = '  -0.51%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.00'
$ws.Range('D15').ClearFormats()
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000172'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -9.98%  '
$ws.Range('D17').Value = '63.837.91'
$ws.Range('E17').Value = '  -1.98%  '
$ws.Range('D18').Value = '3.463.35'
$ws.Range('E18').Value = '  -1.48%  '
$ws.Range('E19').Value = '  -5.52%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.60'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '375.34'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.61%  '
$ws.Range('E22').Value = '  -4.99%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.62'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('E25').Value = '  -7.26%  '
$ws.Range('E26').Value = '  -2.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.64'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.177'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.81%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.97'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -3.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.39'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -7.75%  '
$ws.Range('E32').Value = '  -3.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.79'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -3.48%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.91'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.53'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -5.12%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.34'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.83'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.40%  '
$ws.Range('E38').Value = '  +4.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.29'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0729'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -6.82%  '
$ws.Range('D41').Value = '2.777.41'
$ws.Range('E41').Value = '  -3.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.47'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('E43').Value = '  -8.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.38'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0303'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.28'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -3.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.34'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +7.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '324.43'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('E49').Value = '  -5.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.33'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.821'
$ws.Range('D51').ClearFormats()
